# Update data rows 2-13 (Sending cluster x Target cluster combinations) with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Ednra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.169998666666667
$ws.Range("N2").Value = 9.509996000000001
$ws.Range("O2").Value = 0.06457634599094531
$ws.Range("P2").Value = 0.06457634599094531
$ws.Range("Q2").Value = 538.6324061357907
$ws.Range("R2").Value = 4847.691655222116
$ws.Range("S2").Value = 0.02867954947614604
$ws.Range("T2").Value = 0.02867954947614605

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Ednra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 32.709374
$ws.Range("N3").Value = 98.12812199999999
$ws.Range("O3").Value = 0.6663257858061865
$ws.Range("P3").Value = 0.6663257858061865
$ws.Range("Q3").Value = 5557.834773268718
$ws.Range("R3").Value = 50020.51295941846
$ws.Range("S3").Value = 0.2959276039548591
$ws.Range("T3").Value = 0.2959276039548592

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Ednra"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 169.915657
$ws.Range("H4").Value = 509.746971
$ws.Range("I4").Value = 0.4441184931734509
$ws.Range("J4").Value = 0.4441184931734509
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 13.20978866666667
$ws.Range("N4").Value = 39.629366
$ws.Range("O4").Value = 0.2690978682028682
$ws.Range("P4").Value = 0.2690978682028682
$ws.Range("Q4").Value = 2244.549920127821
$ws.Range("R4").Value = 20200.94928115038
$ws.Range("S4").Value = 0.1195113397424457
$ws.Range("T4").Value = 0.1195113397424457

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Ednra"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 68.382243
$ws.Range("H5").Value = 205.146729
$ws.Range("I5").Value = 0.1787346690539575
$ws.Range("J5").Value = 0.1787346690539575
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.169998666666667
$ws.Range("N5").Value = 9.509996000000001
$ws.Range("O5").Value = 0.06457634599094531
$ws.Range("P5").Value = 0.06457634599094531
$ws.Range("Q5").Value = 216.771619133676
$ws.Range("R5").Value = 1950.944572203084
$ws.Range("S5").Value = 0.01154203182940547
$ws.Range("T5").Value = 0.01154203182940547

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Ednra"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 68.382243
$ws.Range("H6").Value = 205.146729
$ws.Range("I6").Value = 0.1787346690539575
$ws.Range("J6").Value = 0.1787346690539575
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 32.709374
$ws.Range("N6").Value = 98.12812199999999
$ws.Range("O6").Value = 0.6663257858061865
$ws.Range("P6").Value = 0.6663257858061865
$ws.Range("Q6").Value = 2236.740361245882
$ws.Range("R6").Value = 20130.66325121294
$ws.Range("S6").Value = 0.1190955188081869
$ws.Range("T6").Value = 0.1190955188081869

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Ednra"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 68.382243
$ws.Range("H7").Value = 205.146729
$ws.Range("I7").Value = 0.1787346690539575
$ws.Range("J7").Value = 0.1787346690539575
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.20978866666667
$ws.Range("N7").Value = 39.629366
$ws.Range("O7").Value = 0.2690978682028682
$ws.Range("P7").Value = 0.2690978682028682
$ws.Range("Q7").Value = 903.3149785826461
$ws.Range("R7").Value = 8129.834807243813
$ws.Range("S7").Value = 0.04809711841636512
$ws.Range("T7").Value = 0.04809711841636512

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Ednra"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 53.27463399999999
$ws.Range("H8").Value = 159.823902
$ws.Range("I8").Value = 0.1392470275793777
$ws.Range("J8").Value = 0.1392470275793778
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.169998666666667
$ws.Range("N8").Value = 9.509996000000001
$ws.Range("O8").Value = 0.06457634599094531
$ws.Range("P8").Value = 0.06457634599094531
$ws.Range("Q8").Value = 168.8805187471546
$ws.Range("R8").Value = 1519.924668724392
$ws.Range("S8").Value = 0.0089920642311766
$ws.Range("T8").Value = 0.008992064231176601

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Ednra"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 53.27463399999999
$ws.Range("H9").Value = 159.823902
$ws.Range("I9").Value = 0.1392470275793777
$ws.Range("J9").Value = 0.1392470275793778
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 32.709374
$ws.Range("N9").Value = 98.12812199999999
$ws.Range("O9").Value = 0.6663257858061865
$ws.Range("P9").Value = 0.6663257858061865
$ws.Range("Q9").Value = 1742.579928219116
$ws.Range("R9").Value = 15683.21935397204
$ws.Range("S9").Value = 0.0927838850730046
$ws.Range("T9").Value = 0.09278388507300461

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Gnai2"
$ws.Range("C10").Value = "Ednra"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 53.27463399999999
$ws.Range("H10").Value = 159.823902
$ws.Range("I10").Value = 0.1392470275793777
$ws.Range("J10").Value = 0.1392470275793778
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.20978866666667
$ws.Range("N10").Value = 39.629366
$ws.Range("O10").Value = 0.2690978682028682
$ws.Range("P10").Value = 0.2690978682028682
$ws.Range("Q10").Value = 703.7466564340145
$ws.Range("R10").Value = 6333.71990790613
$ws.Range("S10").Value = 0.03747107827519654
$ws.Range("T10").Value = 0.03747107827519655

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Gnai2"
$ws.Range("C11").Value = "Ednra"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 91.01828266666666
$ws.Range("H11").Value = 273.054848
$ws.Range("I11").Value = 0.2378998101932138
$ws.Range("J11").Value = 0.2378998101932138
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.169998666666667
$ws.Range("N11").Value = 9.509996000000001
$ws.Range("O11").Value = 0.06457634599094531
$ws.Range("P11").Value = 0.06457634599094531
$ws.Range("Q11").Value = 288.5278346956231
$ws.Range("R11").Value = 2596.750512260608
$ws.Range("S11").Value = 0.01536270045421719
$ws.Range("T11").Value = 0.01536270045421719

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Gnai2"
$ws.Range("C12").Value = "Ednra"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 91.01828266666666
$ws.Range("H12").Value = 273.054848
$ws.Range("I12").Value = 0.2378998101932138
$ws.Range("J12").Value = 0.2378998101932138
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 32.709374
$ws.Range("N12").Value = 98.12812199999999
$ws.Range("O12").Value = 0.6663257858061865
$ws.Range("P12").Value = 0.6663257858061865
$ws.Range("Q12").Value = 2977.151048581717
$ws.Range("R12").Value = 26794.35943723545
$ws.Range("S12").Value = 0.1585187779701358
$ws.Range("T12").Value = 0.1585187779701358

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Gnai2"
$ws.Range("C13").Value = "Ednra"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 91.01828266666666
$ws.Range("H13").Value = 273.054848
$ws.Range("I13").Value = 0.2378998101932138
$ws.Range("J13").Value = 0.2378998101932138
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 13.20978866666667
$ws.Range("N13").Value = 39.629366
$ws.Range("O13").Value = 0.2690978682028682
$ws.Range("P13").Value = 0.2690978682028682
$ws.Range("Q13").Value = 1202.332278829596
$ws.Range("R13").Value = 10820.99050946637
$ws.Range("S13").Value = 0.06401833176886081
$ws.Range("T13").Value = 0.06401833176886081

# Remove the now-obsolete rows (old "Resolving-Mac" sending-cluster block, rows 14-17)
$ws.Range("A14:T17").Delete()
